$wb = $excel.ActiveWorkbook

# Update "Latest Handback DateTime" (column K) for the 30c1937e row (row 2)
# on both the zh-cn and de-de localization-status sheets, reflecting a
# freshly generated handback report.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-09-07 03:40:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-09-07 03:40:54"
